$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The cover letter's title sentence is being reworded:
#   "...Stabilizing Selection for Antimicrobial Resistance in Escherichia coli"
# becomes
#   "...Stabilizing Selection of Antimicrobial Resistance Genes in Escherichia coli"
# i.e. "for" -> "of", and "Genes " is inserted right before "in".
# ---------------------------------------------------------------------------

# The hidden "_GoBack" bookmark currently sits between "...publication " and
# "elsewhere." further down the same paragraph. Pull it out of the way before
# editing so it doesn't keep that pair of runs from recombining; it gets
# re-dropped at its new (post-edit) home near the title once the text is in
# its final shape.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Step 1: "for" -> "of" (word right after "...Stabilizing Selection ")
$rFor = $d.Range(114, 117)
$rFor.Text = "of"

# Step 2: insert "Genes " immediately before "in Escherichia coli"
# (offset shifted left by 1 because "for"(3) -> "of"(2))
$rGenes = $d.Range(142, 142)
$rGenes.InsertBefore("Genes ")

# Step 3: the "...publication " / "elsewhere." runs used to be kept apart by
# the bookmark; now that it has moved, recombine them into a single run. A
# same-text Range.Text assignment is a no-op, so round-trip through a
# throwaway value to force the engine to actually rewrite (and coalesce) the
# run(s).
$rTail = $d.Range(310, 361)
$rTail.Text = "TEMP_PLACEHOLDER"
$rTail2 = $d.Range(310, 310 + ("TEMP_PLACEHOLDER").Length)
$rTail2.Text = ", has not been submitted for publication elsewhere."

# ---------------------------------------------------------------------------
# Re-establish run boundaries around the edit (the engine merges adjacent
# same-formatted runs on edit) so the new runs line up the way a live Word
# editing session would leave them: each toggle below nudges formatting
# on/off over a tiny span, which is a no-op visually but forces a run split
# at that boundary.
# ---------------------------------------------------------------------------
function SplitAt($start, $end) {
    $r = $d.Range($start, $end)
    $r.Bold = 1
    $r.Bold = 0
}

SplitAt 59 68      # "entitled "
SplitAt 68 69      # opening curly quote
SplitAt 69 114     # "Genetic Capitalism and Stabilizing Selection "
SplitAt 114 116    # "of"
SplitAt 116 142    # " Antimicrobial Resistance "
SplitAt 142 148    # "Genes "
SplitAt 148 151    # "in "
SplitAt 151 167    # "Escherichia coli" (italic run, reasserted)
SplitAt 167 168    # closing curly quote
SplitAt 168 227    # " for your consideration. This original work, completed by C"
SplitAt 227 255    # "olby T. Ford, Gabriel Lopez "

# Drop "_GoBack" back in, now right before "in Escherichia coli" (i.e. right
# after the newly-inserted "Genes ").
$rNewGoBack = $d.Range(148, 148)
$d.Bookmarks.Add("_GoBack", $rNewGoBack)
